$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 70, shifting existing rows 70:278 down to 71:279
$ws.Rows.Item(70).Insert()

# Populate the newly inserted row 70 with the new weekly record
$ws.Cells.Item(70, 1).Value = 3
$ws.Cells.Item(70, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(70, 3).Value = "Coquimbo"
$ws.Cells.Item(70, 4).Value = 44624
$ws.Cells.Item(70, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(70, 5).Value = 5
$ws.Cells.Item(70, 6).Value = 100112039
$ws.Cells.Item(70, 7).Value = "Ciboulette"
$ws.Cells.Item(70, 8).Value = "Sin especificar"
$ws.Cells.Item(70, 9).Value = "Primera"
$ws.Cells.Item(70, 10).Value = 160
$ws.Cells.Item(70, 11).Value = 1500
$ws.Cells.Item(70, 12).Value = 1500
$ws.Cells.Item(70, 13).Value = 1500
$ws.Cells.Item(70, 14).Value = "$/docena de atados"
$ws.Cells.Item(70, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(70, 16).Value = 500
$ws.Cells.Item(70, 17).Value = 3
$ws.Cells.Item(70, 18).Value = "Hortaliza"
